$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tester")
$ws.Range("D4").NumberFormat = "0.00000"
